{"js": "// The commit inserts \" of the game's KickStarter\" right after the\n// phrase \" (only backers\" (and before \" have received PDF copies...\")\n// in the Introduction paragraph.\n\nconst body = context.document.body;\n\n// Find the unique anchor phrase so the insertion point is precise.\nconst results = body.search(\" (only backers\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Anchor phrase ' (only backers' not found in document body.\");\n}\n\n// Insert the new text immediately after the matched phrase, pushing the\n// rest of the sentence (\" have received PDF copies...\") further along.\nconst target = results.items[0];\ntarget.insertText(\" of the game\\u2019s KickStarter\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# The commit inserts \" of the game's KickStarter\" right after the\n# phrase \" (only backers\" (and before \" have received PDF copies...\")\n# in the Introduction paragraph.\n\n$d = $word.ActiveDocument\n$range = $d.Content\n\n$findText = \" (only backers\"\n$replaceText = \" (only backers of the game\" + [char]0x2019 + \"s KickStarter\"\n\n$found = $range.Find.Execute(\n    $findText,    # FindText\n    $true,        # MatchCase\n    $false,       # MatchWholeWord\n    $false,       # MatchWildcards\n    $false,       # MatchSoundsLike\n    $false,       # MatchAllWordForms\n    $true,        # Forward\n    1,            # Wrap (wdFindContinue)\n    $false,       # Format\n    $replaceText, # ReplaceWith\n    1             # Replace (wdReplaceOne)\n)\n\nif (-not $found) {\n    throw \"Anchor phrase ' (only backers' not found in document body.\"\n}\n"}
